$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '68.753.83'
$ws.Range('E2').Value = '  +1.78%  '
Set-TextValue $ws.Range('D3') '3.279.26'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue $ws.Range('D5') '584.64'
$ws.Range('E5').Value = '  +0.99%  '
Set-TextValue $ws.Range('D6') '182.55'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue $ws.Range('D8') '0.598'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('E9').Value = '  +2.98%  '
Set-TextValue $ws.Range('D10') '6.65'
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  +1.82%  '
Set-TextValue $ws.Range('D12') '3.860.70'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('E13').Value = '  -0.32%  '
Set-TextValue $ws.Range('D14') '28.83'
$ws.Range('E14').Value = '  +1.30%  '
Set-TextValue $ws.Range('D15') '68.719.45'
$ws.Range('E15').Value = '  +1.76%  '
Set-TextValue $ws.Range('D16') '0.0000171'
$ws.Range('E16').Value = '  +2.51%  '
Set-TextValue $ws.Range('D17') '3.270.56'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('E18').Value = '  +0.22%  '
Set-TextValue $ws.Range('D19') '13.57'
$ws.Range('E19').Value = '  +0.39%  '
Set-TextValue $ws.Range('D20') '394.29'
$ws.Range('E20').Value = '  +4.62%  '
$ws.Range('E21').Value = '  +1.41%  '
Set-TextValue $ws.Range('D22') '71.49'
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue $ws.Range('D23') '0.998'
$ws.Range('E23').Value = '  -0.12%  '
Set-TextValue $ws.Range('D24') '0.515'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('E26').Value = '  +3.84%  '
Set-TextValue $ws.Range('D27') '9.64'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('E28').Value = '  -0.68%  '
Set-TextValue $ws.Range('D29') '1.99'
$ws.Range('E29').Value = '  +0.84%  '
Set-TextValue $ws.Range('D30') '5.73'
$ws.Range('E30').Value = '  +0.21%  '
Set-TextValue $ws.Range('D31') '22.99'
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D32') '1.30'
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D33') '7.14'
$ws.Range('E33').Value = '  +3.42%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  +1.74%  '
Set-TextValue $ws.Range('D36') '164.01'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('E37').Value = '  +1.10%  '
Set-TextValue $ws.Range('D38') '0.832'
$ws.Range('E38').Value = '  -2.39%  '
Set-TextValue $ws.Range('D39') '4.59'
$ws.Range('E39').Value = '  +2.21%  '
Set-TextValue $ws.Range('D40') '26.10'
$ws.Range('E40').Value = '  -2.51%  '
Set-TextValue $ws.Range('D41') '6.56'
$ws.Range('E41').Value = '  -3.19%  '
Set-TextValue $ws.Range('D42') '2.55'
$ws.Range('E42').Value = '  -2.84%  '
Set-TextValue $ws.Range('D43') '41.49'
$ws.Range('E43').Value = '  +2.15%  '
Set-TextValue $ws.Range('D44') '0.0687'
$ws.Range('E44').Value = '  +1.84%  '
Set-TextValue $ws.Range('D45') '343.16'
Set-TextValue $ws.Range('D46') '2.606.04'
$ws.Range('E46').Value = '  -4.99%  '
Set-TextValue $ws.Range('D47') '24.76'
$ws.Range('E47').Value = '  -2.81%  '
Set-TextValue $ws.Range('D48') '0.0281'
$ws.Range('E48').Value = '  +0.96%  '
Set-TextValue $ws.Range('D49') '31.97'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('E51').Value = '  -0.33%  '
